$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.589.18"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.597.34"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.09"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.65"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  +2.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.66"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.049.84"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "60.552.72"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.74"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.596.43"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.75"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "351.99"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.60"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.23"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.02"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.429"
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.712.80"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.23"
$ws.Range("E31").Value = "  +8.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.44"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.83"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.11"
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.917"
$ws.Range("E37").Value = "  +6.63%  "
$ws.Range("E38").Value = "  +2.88%  "
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.37"
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.841"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.70"
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.59"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.82"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.974.44"
$ws.Range("E51").Value = "  -1.88%  "
